{"js": "// Apply the redline/placeholder-fill edits described by the diff:\n//   - Fill in the letterhead name/address/city/state/zip, SSN, and DOB.\n//   - Fill in the account/dispute-reason bullet placeholder.\n//   - Fill in the signature name and trailing SSN.\n//\n// Each replacement is done via a literal (non-wildcard, case-sensitive)\n// body.search() for the exact old text, then insertText(..., \"Replace\") on\n// the single hit. The SSN-line replacement runs before the standalone\n// \"[ss_number]\" replacement so that the later search is unambiguous (the\n// placeholder appears twice in the original document).\n\nasync function replaceOnce(body, oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for ${JSON.stringify(oldText)}, found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\nconst body = context.document.body;\n\nawait replaceOnce(body, \"John [client_middle_name] Doe\", \"sandesh nothing sitaula\");\nawait replaceOnce(body, \"[client_address]\", \"Surunga kanakai2\");\nawait replaceOnce(body, \"[client_city], [client_state] [client_postal_code]\", \"Kanakai, Koshi 50354\");\nawait replaceOnce(body, \"SSN: [ss_number] DOB: [bdate]\", \"SSN: 01010101 DOB: 04-12-2024\");\nawait replaceOnce(body, \"[account , dispute_reason_in_bullet_list]\", \"some accoutn detail there was som edisput\");\nawait replaceOnce(body, \"John Doe\", \"sandesh sitaula\");\nawait replaceOnce(body, \"[ss_number]\", \"01010101\");\n", "ps1": "# Apply the redline/placeholder-fill edits described by the diff:\n#   - Fill in the letterhead name/address/city/state/zip, SSN, and DOB.\n#   - Fill in the account/dispute-reason bullet placeholder.\n#   - Fill in the signature name and trailing SSN.\n#\n# Each replacement runs Find/Replace (literal text, case-sensitive, no\n# wildcards) over the whole document body. The SSN-line replacement runs\n# before the standalone \"[ss_number]\" replacement so that by the time the\n# second Find runs, \"[ss_number]\" is unique again (it appears twice in the\n# original document).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #          MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n    #          Format, ReplaceWith, Replace)\n    # Wrap:=1 (wdFindContinue), Replace:=2 (wdReplaceAll)\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-Text \"John [client_middle_name] Doe\" \"sandesh nothing sitaula\"\nReplace-Text \"[client_address]\" \"Surunga kanakai2\"\nReplace-Text \"[client_city], [client_state] [client_postal_code]\" \"Kanakai, Koshi 50354\"\nReplace-Text \"SSN: [ss_number] DOB: [bdate]\" \"SSN: 01010101 DOB: 04-12-2024\"\nReplace-Text \"[account , dispute_reason_in_bullet_list]\" \"some accoutn detail there was som edisput\"\nReplace-Text \"John Doe\" \"sandesh sitaula\"\nReplace-Text \"[ss_number]\" \"01010101\"\n"}
